# Generate Report for Handoff
# Rename the handed-off source file from
#   fa3ccb9a-d0e6-4c60-be32-4eb6a3ac83ec.md
# to
#   41c0a2fa-7ee3-4a4e-8dae-a6c0d4ac5d8f.md
# across the Overview / zh-cn / de-de sheets, bump the recorded
# handoff/handback timestamps, and rename the generated .xlf artifacts
# to carry the new content hash.

$wb = $excel.ActiveWorkbook

$oldBase = "fa3ccb9a-d0e6-4c60-be32-4eb6a3ac83ec"
$newBase = "41c0a2fa-7ee3-4a4e-8dae-a6c0d4ac5d8f"

$oldMd = $oldBase + ".md"
$newMd = $newBase + ".md"

$ghPrefix = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d914958733067bab925aef64d616dd7e1d22d5b0/e2e/"
$newGhUrl = $ghPrefix + $newMd

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newMd
$wsOverview.Range("B2").Value = "e2e\" + $newMd
$wsOverview.Range("G2").Value = "2016-09-02 11:08:02"

# Update the hyperlink on B2 in place (keeps the cell's existing
# "HyperLink" formatting untouched) instead of re-adding a new link.
$wsOverview.Hyperlinks.Delete()
$hlOverview = $wsOverview.Range("B2").Hyperlinks.Item(1)
$hlOverview.Address = $newGhUrl
$hlOverview.TextToDisplay = "e2e\" + $newMd

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = $newMd
$wsZhCn.Range("G2").Value = $newBase + ".c4004000eca71d9d6295042713ef2b07d8932e65.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-09-02 11:07:57"

$wsZhCn.Hyperlinks.Delete()
$hlZhCn = $wsZhCn.Range("A2").Hyperlinks.Item(1)
$hlZhCn.Address = $newGhUrl
$hlZhCn.TextToDisplay = $newMd

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = $newMd
$wsDeDe.Range("G2").Value = $newBase + ".c4004000eca71d9d6295042713ef2b07d8932e65.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-09-02 11:08:02"

$wsDeDe.Hyperlinks.Delete()
$hlDeDe = $wsDeDe.Range("A2").Hyperlinks.Item(1)
$hlDeDe.Address = $newGhUrl
$hlDeDe.TextToDisplay = $newMd
